# Add two new Mac-Addresses (10 new rows) to the reg_center_machine_device sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for the new rows: regcntr_id, machine_id, device_id, lang_code, is_active, cr_by, cr_dtimes
$newRows = @(
    @(10001, 10030, 3000166, "eng", $true, "superadmin", "now()"),
    @(10001, 10030, 3000167, "eng", $true, "superadmin", "now()"),
    @(10001, 10030, 3000168, "eng", $true, "superadmin", "now()"),
    @(10001, 10030, 3000169, "eng", $true, "superadmin", "now()"),
    @(10001, 10030, 3000170, "eng", $true, "superadmin", "now()"),
    @(10001, 10031, 3000171, "eng", $true, "superadmin", "now()"),
    @(10001, 10031, 3000172, "eng", $true, "superadmin", "now()"),
    @(10001, 10031, 3000173, "eng", $true, "superadmin", "now()"),
    @(10001, 10031, 3000174, "eng", $true, "superadmin", "now()"),
    @(10001, 10031, 3000175, "eng", $true, "superadmin", "now()")
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Update the view to reflect the new scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 144
$ws.Range("H149").Select()
